$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L ("Open") rows 2-7 need to switch from the text "true" to the
# text "false" (these are literal text values sharing the same styling as
# the existing "false" cells in rows 11-13, not Excel Booleans).
#
# A plain `$range.Value = "false"` assignment gets auto-coerced by Excel
# into the Boolean FALSE (t="b"), same as typing FALSE straight into a
# cell formatted as Text. The reliable way to land literal text without
# touching the cell's style/number format is: compute the text via a
# formula, copy that computed result, then Paste Special > Values only
# into the target cells.
$helper = $ws.Range("N1")
$helper.Formula = "=""false"""
$helper.Copy()
foreach ($r in 2..7) {
    $ws.Range("L$r").PasteSpecial($excel.Constants.xlPasteValues)
}
$helper.Clear()

# The sheet's saved selection/active cell moves from I16 to K17.
$ws.Range("K17").Select()
